$d = $word.ActiveDocument

$d.Content.Find.Execute("New guide!", $true, $false, $false, $false, $false,
                         $true, 1, $false, "New guide!", 2)

$d.Content.Find.Execute("Tom Coleman", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Tom Coleman", 2)

$d.Content.Find.Execute("Guide on rationalizing the denominator available now!", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Guide on rationalizing the denominator available now!", 2)
